$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new parameter row (row 27)
$ws.Range("A27").Value = "debut_porte_couteaux (mm)"
$ws.Range("B27").Value = 10

# Move active selection to B28, matching the author's last cursor position
$ws.Range("B28").Select()
